$d = $word.ActiveDocument

# The first paragraph currently reads "gmkmnlkgfmbkgfmbkgmblkdfmlk".
# Split it into two runs: "G" and "mkmnlkgfmbkgfmbkgmblkdfmlk" (capitalize first letter).
$para1 = $d.Paragraphs(1)
$r = $para1.Range
$r.Text = "G"

$endOfPara1 = $para1.Range.End - 1
$insertRange = $d.Range($endOfPara1, $endOfPara1)
$insertRange.InsertAfter("mkmnlkgfmbkgfmbkgmblkdfmlk")

# Add a new paragraph after the first one with the new text.
$endOfPara1Full = $d.Paragraphs(1).Range.End
$newParaRange = $d.Range($endOfPara1Full, $endOfPara1Full)
$newParaRange.InsertParagraphAfter()
$newParaRange2 = $d.Paragraphs(2).Range
$newParaRange2.InsertAfter("eighjogfjiofdjgoidfjgoidsjfoisdjfoisj")
